$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44574
$ws.Range("K2").Value = 'Black Amber'
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 19000
$ws.Range("P2").Value = 18500
$ws.Range("Q2").Value = '$/bandeja 18 kilos granel'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 1028

$ws.Range("D3").Value = 44245
$ws.Range("K3").Value = 'Black Amber'
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("Q3").Value = '$/bandeja 18 kilos granel'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 806

$ws.Range("D4").Value = 44174
$ws.Range("K4").Value = 'Angeleno'
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 270
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20500
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 1139

$ws.Range("D5").Value = 44175
$ws.Range("K5").Value = 'Angeleno'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("Q5").Value = '$/bandeja 18 kilos granel'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 1194

$ws.Range("D6").Value = 44285
$ws.Range("K6").Value = 'Angeleno'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14500
$ws.Range("Q6").Value = '$/bandeja 18 kilos granel'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 806

$ws.Range("D7").Value = 44614
$ws.Range("K7").Value = 'Angeleno'
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 19000
$ws.Range("P7").Value = 18500
$ws.Range("Q7").Value = '$/bandeja 18 kilos granel'
$ws.Range("R7").Value = 'Región Metropolitana'
$ws.Range("S7").Value = 1028

$ws.Range("D9").Value = 44278
$ws.Range("K9").Value = 'Angeleno'
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 15500
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 861

$ws.Range("D10").Value = 44596
$ws.Range("K10").Value = 'Black Amber'
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15500
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 861

$ws.Range("D11").Value = 44580
$ws.Range("K11").Value = 'Black Amber'
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 270
$ws.Range("N11").Value = 19000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 19500
$ws.Range("Q11").Value = '$/bandeja 18 kilos granel'
$ws.Range("R11").Value = 'Región Metropolitana'
$ws.Range("S11").Value = 1083

$ws.Range("D12").Value = 44169
$ws.Range("K12").Value = 'Angeleno'
$ws.Range("L12").Value = 'Tercera'
$ws.Range("M12").Value = 250
$ws.Range("N12").Value = 24000
$ws.Range("O12").Value = 25000
$ws.Range("P12").Value = 24500
$ws.Range("Q12").Value = '$/bandeja 18 kilos granel'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 1361

$ws.Range("D13").Value = 44628
$ws.Range("K13").Value = 'Black Amber'
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 270
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 16000
$ws.Range("P13").Value = 15500
$ws.Range("Q13").Value = '$/bandeja 18 kilos granel'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 861

$ws.Range("D14").Value = 44238
$ws.Range("K14").Value = 'Black Amber'
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 14000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 14500
$ws.Range("Q14").Value = '$/bandeja 18 kilos granel'
$ws.Range("R14").Value = 'Región de O''Higgins'
$ws.Range("S14").Value = 806

$ws.Range("D15").Value = 44238
$ws.Range("K15").Value = 'Fortuna'
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 14000
$ws.Range("O15").Value = 15000
$ws.Range("P15").Value = 14500
$ws.Range("Q15").Value = '$/bandeja 18 kilos granel'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 806

$ws.Range("D16").Value = 44314
$ws.Range("K16").Value = 'Angeleno'
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 250
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 14500
$ws.Range("Q16").Value = '$/bandeja 18 kilos granel'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 806

$ws.Range("D17").Value = 44217
$ws.Range("K17").Value = 'Black Amber'
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 300
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 17000
$ws.Range("P17").Value = 16500
$ws.Range("Q17").Value = '$/bandeja 18 kilos granel'
$ws.Range("R17").Value = 'Región Metropolitana'
$ws.Range("S17").Value = 917

$ws.Range("D18").Value = 44650
$ws.Range("K18").Value = 'Angeleno'
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = 17000
$ws.Range("O18").Value = 18000
$ws.Range("P18").Value = 17500
$ws.Range("Q18").Value = '$/bandeja 18 kilos granel'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 972

$ws.Range("D19").Value = 44239
$ws.Range("K19").Value = 'Fortuna'
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 300
$ws.Range("N19").Value = 15000
$ws.Range("O19").Value = 16000
$ws.Range("P19").Value = 15500
$ws.Range("Q19").Value = '$/bandeja 18 kilos granel'
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 861

$ws.Range("D20").Value = 44243
$ws.Range("K20").Value = 'Black Amber'
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 300
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 15000
$ws.Range("P20").Value = 14500
$ws.Range("Q20").Value = '$/caja 18 kilos granel'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 806

$ws.Range("D21").Value = 44587
$ws.Range("K21").Value = 'Black Amber'
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 300
$ws.Range("N21").Value = 15000
$ws.Range("O21").Value = 16000
$ws.Range("P21").Value = 15500
$ws.Range("Q21").Value = '$/caja 18 kilos granel'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 861
